$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "91-62=29"
$t.Cell(1,2).Range.Text = "77-11=66"
$t.Cell(1,3).Range.Text = "69-33=36"
$t.Cell(1,4).Range.Text = "74-31=43"
$t.Cell(1,5).Range.Text = "1+43=44"
$t.Cell(2,1).Range.Text = "39+0=39"
$t.Cell(2,2).Range.Text = "50+23=73"
$t.Cell(2,3).Range.Text = "78-76=2"
$t.Cell(2,4).Range.Text = "14+69=83"
$t.Cell(2,5).Range.Text = "0+35=35"
$t.Cell(3,1).Range.Text = "88-70=18"
$t.Cell(3,2).Range.Text = "73-26=47"
$t.Cell(3,3).Range.Text = "68-4=64"
$t.Cell(3,4).Range.Text = "76-36=40"
$t.Cell(3,5).Range.Text = "40+33=73"
$t.Cell(4,1).Range.Text = "42+18=60"
$t.Cell(4,2).Range.Text = "26+13=39"
$t.Cell(4,3).Range.Text = "14+5=19"
$t.Cell(4,4).Range.Text = "25+15=40"
$t.Cell(4,5).Range.Text = "93-40=53"
$t.Cell(5,1).Range.Text = "83-62=21"
$t.Cell(5,2).Range.Text = "28+12=40"
$t.Cell(5,3).Range.Text = "64-23=41"
$t.Cell(5,4).Range.Text = "2+63=65"
$t.Cell(5,5).Range.Text = "56+10=66"
$t.Cell(6,1).Range.Text = "57-40=17"
$t.Cell(6,2).Range.Text = "69+7=76"
$t.Cell(6,3).Range.Text = "45+26=71"
$t.Cell(6,4).Range.Text = "12-8=4"
$t.Cell(6,5).Range.Text = "48+36=84"
$t.Cell(7,1).Range.Text = "50+12=62"
$t.Cell(7,2).Range.Text = "46+4=50"
$t.Cell(7,3).Range.Text = "77+9=86"
$t.Cell(7,4).Range.Text = "31+20=51"
$t.Cell(7,5).Range.Text = "14+52=66"
$t.Cell(8,1).Range.Text = "93-59=34"
$t.Cell(8,2).Range.Text = "12+45=57"
$t.Cell(8,3).Range.Text = "33+47=80"
$t.Cell(8,4).Range.Text = "53-52=1"
$t.Cell(8,5).Range.Text = "41-9=32"
$t.Cell(9,1).Range.Text = "17-17=0"
$t.Cell(9,2).Range.Text = "76-40=36"
$t.Cell(9,3).Range.Text = "18+65=83"
$t.Cell(9,4).Range.Text = "93+4=97"
$t.Cell(9,5).Range.Text = "3+29=32"
$t.Cell(10,1).Range.Text = "5+88=93"
$t.Cell(10,2).Range.Text = "93-24=69"
$t.Cell(10,3).Range.Text = "41-28=13"
$t.Cell(10,4).Range.Text = "47-31=16"
$t.Cell(10,5).Range.Text = "86-4=82"
$t.Cell(11,1).Range.Text = "31-23=8"
$t.Cell(11,2).Range.Text = "6+37=43"
$t.Cell(11,3).Range.Text = "90-43=47"
$t.Cell(11,4).Range.Text = "10-6=4"
$t.Cell(11,5).Range.Text = "3+66=69"
$t.Cell(12,1).Range.Text = "27+7=34"
$t.Cell(12,2).Range.Text = "36+45=81"
$t.Cell(12,3).Range.Text = "13+40=53"
$t.Cell(12,4).Range.Text = "19-6=13"
$t.Cell(12,5).Range.Text = "22+48=70"
$t.Cell(13,1).Range.Text = "51+48=99"
$t.Cell(13,2).Range.Text = "14+1=15"
$t.Cell(13,3).Range.Text = "4+76=80"
$t.Cell(13,4).Range.Text = "63-3=60"
$t.Cell(13,5).Range.Text = "94-50=44"
$t.Cell(14,1).Range.Text = "42-18=24"
$t.Cell(14,2).Range.Text = "12+31=43"
$t.Cell(14,3).Range.Text = "3+27=30"
$t.Cell(14,4).Range.Text = "81-4=77"
$t.Cell(14,5).Range.Text = "20+48=68"
$t.Cell(15,1).Range.Text = "68+5=73"
$t.Cell(15,2).Range.Text = "81+17=98"
$t.Cell(15,3).Range.Text = "71+10=81"
$t.Cell(15,4).Range.Text = "7+89=96"
$t.Cell(15,5).Range.Text = "41+27=68"
$t.Cell(16,1).Range.Text = "17+33=50"
$t.Cell(16,2).Range.Text = "8+8=16"
$t.Cell(16,3).Range.Text = "5+51=56"
$t.Cell(16,4).Range.Text = "23-14=9"
$t.Cell(16,5).Range.Text = "66-58=8"
$t.Cell(17,1).Range.Text = "18+28=46"
$t.Cell(17,2).Range.Text = "53+3=56"
$t.Cell(17,3).Range.Text = "83-13=70"
$t.Cell(17,4).Range.Text = "27+46=73"
$t.Cell(17,5).Range.Text = "18+60=78"
$t.Cell(18,1).Range.Text = "17+63=80"
$t.Cell(18,2).Range.Text = "60+22=82"
$t.Cell(18,3).Range.Text = "79-39=40"
$t.Cell(18,4).Range.Text = "37+34=71"
$t.Cell(18,5).Range.Text = "21-20=1"
$t.Cell(19,1).Range.Text = "9+85=94"
$t.Cell(19,2).Range.Text = "95-31=64"
$t.Cell(19,3).Range.Text = "24-5=19"
$t.Cell(19,4).Range.Text = "87-39=48"
$t.Cell(19,5).Range.Text = "85-41=44"
$t.Cell(20,1).Range.Text = "27+45=72"
$t.Cell(20,2).Range.Text = "91-48=43"
$t.Cell(20,3).Range.Text = "85-30=55"
$t.Cell(20,4).Range.Text = "85+0=85"
$t.Cell(20,5).Range.Text = "25+70=95"
